$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store numeric-looking / percent-looking
# text as plain strings (t="inlineStr") in the original workbook. Excel's COM
# layer auto-converts numeric-looking strings assigned via .Value into real
# numbers, so we temporarily force the whole D2:E51 block to Text format,
# write the new string values, then restore the default ("Normal") style so
# no stray style index is left behind on the cells.
$numRng = $ws.Range("D2:E51")
$numRng.NumberFormat = "@"

$ws.Range("D2").Value = "308.08"
$ws.Range("E2").Value = "1.98%"

$ws.Range("D3").Value = "36.22"
$ws.Range("E3").Value = "3.29%"

$ws.Range("D4").Value = "5.101"
$ws.Range("E4").Value = "1.12%"

$ws.Range("D5").Value = "0.08127"
$ws.Range("E5").Value = "2.81%"

$ws.Range("D6").Value = "1.942"
$ws.Range("E6").Value = "1.59%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.187"
$ws.Range("E7").Value = "3.93%"

$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "7.794"
$ws.Range("E8").Value = "1.03%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9312"
$ws.Range("E9").Value = "0.90%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1393"
$ws.Range("E10").Value = "16.57%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1926"
$ws.Range("E11").Value = "4.83%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09204"
$ws.Range("E12").Value = "-2.73%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03392"
$ws.Range("E13").Value = "-4.05%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09864"
$ws.Range("E14").Value = "-0.21%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001408"
$ws.Range("E15").Value = "0.52%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005804"
$ws.Range("E16").Value = "-0.30%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.605"
$ws.Range("E17").Value = "3.23%"

$ws.Range("D18").Value = "2.986"
$ws.Range("E18").Value = "-0.03%"

$ws.Range("D19").Value = "0.3440"
$ws.Range("E19").Value = "-0.08%"

$ws.Range("D20").Value = "0.1350"
$ws.Range("E20").Value = "4.60%"

$ws.Range("D21").Value = "4.886"
$ws.Range("E21").Value = "-3.03%"

$ws.Range("D23").Value = "0.04516"
$ws.Range("E23").Value = "0.45%"

$ws.Range("D24").Value = "0.001217"
$ws.Range("E24").Value = "-0.02%"

$ws.Range("D25").Value = "0.004863"
$ws.Range("E25").Value = "6.51%"

$ws.Range("D26").Value = "0.0001240"
$ws.Range("E26").Value = "-0.73%"

$ws.Range("D40").Value = "0.04950"
$ws.Range("E40").Value = "4.90%"

$ws.Range("D41").Value = "0.007639"
$ws.Range("E41").Value = "0.51%"

$ws.Range("D42").Value = "0.01026"
$ws.Range("E42").Value = "7.40%"

$ws.Range("D43").Value = "0.1384"
$ws.Range("E43").Value = "4.54%"

$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").Value = "-0.41%"

$ws.Range("D45").Value = "0.01135"
$ws.Range("E45").Value = "1.52%"

$ws.Range("D46").Value = "0.00006444"
$ws.Range("E46").Value = "4.39%"

$ws.Range("E47").Value = "0.14%"

$ws.Range("E49").Value = "-8.60%"

$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "0.14%"

$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").Value = "0.14%"

# Restore default styling on the numeric block so no explicit style index
# is left on any cell (matches the original unstyled inlineStr cells).
$numRng.Style = "Normal"
